# Weekly driver report update for 2025-04-28
# Updates the "HarrisHealth_driver_summary" Driver Summary sheet:
#  - Bad Drivers table: new driver rows + recomputed totals
#  - Good Drivers table: refreshed driver rows (23 entries)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlRight = -4152

# ---------------------------------------------------------------------------
# 0) Capture formatting from stable cells that keep their style, BEFORE we
#    overwrite the cells that currently hold styles 5 (bold+#,##0) so we can
#    still stamp that look onto the new "Totals:" row further down.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 2).Copy() | Out-Null          # style 5 (bold, #,##0)
$ws.Cells.Item(5, 2).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(4, 3).Copy() | Out-Null          # style 5 (bold, #,##0)
$ws.Cells.Item(5, 3).PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 1) Bad Drivers table (rows 3-5)
# ---------------------------------------------------------------------------
# Row 3 - existing driver row gets new figures
$ws.Cells.Item(3, 1).Value = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.40.0.4"
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 66.7

# Row 4 - new driver row (style 4: right aligned, no border/bold)
$ws.Cells.Item(4, 1).Value = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.32.1"
$ws.Cells.Item(4, 1).ClearFormats() | Out-Null
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 98.3
$ws.Cells.Item(3, 2).Copy() | Out-Null          # style 4 source
$ws.Cells.Item(4, 2).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(3, 3).Copy() | Out-Null
$ws.Cells.Item(4, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(3, 4).Copy() | Out-Null
$ws.Cells.Item(4, 4).PasteSpecial($xlPasteFormats) | Out-Null

# Row 5 - Totals row (moved down from row 4); format already pasted in step 0
$ws.Cells.Item(5, 1).Value = "Totals:"
$ws.Cells.Item(1, 1).Copy() | Out-Null          # style 1 (bold) source
$ws.Cells.Item(5, 1).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 3
$ws.Cells.Item(5, 4).ClearContents() | Out-Null
$ws.Cells.Item(5, 4).ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 2) Clear the old row-4 bold style and the old row-10 section header so the
#    rows that are now blank (6-10) carry no leftover formatting.
# ---------------------------------------------------------------------------
$ws.Cells.Item(10, 1).ClearContents() | Out-Null
$ws.Cells.Item(10, 1).ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 3) Good Drivers section header moves from row 10 -> row 11
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "Good Drivers (Roaming > 99.8%)"
$ws.Cells.Item(1, 1).Copy() | Out-Null          # style 1 (bold) source
$ws.Cells.Item(11, 1).PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 4) Column header row moves from row 11 -> row 12
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "Adapter-Driver"
$ws.Cells.Item(2, 1).Copy() | Out-Null          # style 2 (bottom border)
$ws.Cells.Item(12, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(12, 2).Value = "Total Samples"
$ws.Cells.Item(2, 2).Copy() | Out-Null          # style 3 (border + right)
$ws.Cells.Item(12, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(12, 3).Value = ""
$ws.Cells.Item(2, 1).Copy() | Out-Null          # style 2 (bottom border)
$ws.Cells.Item(12, 3).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(12, 4).Value = "Good Roaming Calculation (%)"
$ws.Cells.Item(2, 2).Copy() | Out-Null          # style 3 (border + right)
$ws.Cells.Item(12, 4).PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(12, 5).Value = "Driver Vintage"
$ws.Cells.Item(2, 2).Copy() | Out-Null          # style 3 (border + right)
$ws.Cells.Item(12, 5).PasteSpecial($xlPasteFormats) | Out-Null

# Old header row (row 11 content) no longer exists as a header; clear the
# trailing cells that used to hold it (column A was already reused above).
$ws.Cells.Item(11, 2).ClearContents() | Out-Null
$ws.Cells.Item(11, 2).ClearFormats() | Out-Null
$ws.Cells.Item(11, 3).ClearContents() | Out-Null
$ws.Cells.Item(11, 3).ClearFormats() | Out-Null
$ws.Cells.Item(11, 4).ClearContents() | Out-Null
$ws.Cells.Item(11, 4).ClearFormats() | Out-Null
$ws.Cells.Item(11, 5).ClearContents() | Out-Null
$ws.Cells.Item(11, 5).ClearFormats() | Out-Null

# ---------------------------------------------------------------------------
# 5) Good Drivers data rows 13-35 (23 entries). Rows 13-17 already carry the
#    right styles (6 / 4) from the previous data; rows 18-35 are brand new
#    and need formatting copied in from row 13.
# ---------------------------------------------------------------------------
$driverRows = @(
    @{ Row = 13; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 23.70.2.3";          Samples = 18721;  Pct = 99.90000000000001; Vintage = "2024-07-23" }
    @{ Row = 14; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.250.10.1";        Samples = 69578;  Pct = 99.90000000000001; Vintage = "2023-08-14" }
    @{ Row = 15; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.230.0.8";         Samples = 338880; Pct = 99.90000000000001; Vintage = "2023-05-08" }
    @{ Row = 16; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.200.0.6";         Samples = 143869; Pct = 99.90000000000001; Vintage = "2023-01-16" }
    @{ Row = 17; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.190.0.4";         Samples = 287148; Pct = 99.90000000000001; Vintage = "2022-11-22" }
    @{ Row = 18; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.160.0.4";         Samples = 96526;  Pct = 99.90000000000001; Vintage = "2022-08-13" }
    @{ Row = 19; Name = "Intel(R) Dual Band Wireless-AC 8265 - 22.30.0.11";   Samples = 172690; Pct = 99.90000000000001; Vintage = "2021-01-19" }
    @{ Row = 20; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.30.0.11";         Samples = 67111;  Pct = 100;                Vintage = "2021-01-19" }
    @{ Row = 21; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.10.0.7";          Samples = 68450;  Pct = 100;                Vintage = "2020-10-19" }
    @{ Row = 22; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 22.0.1.1";           Samples = 15734;  Pct = 99.90000000000001; Vintage = "2020-09-28" }
    @{ Row = 23; Name = "Intel(R) Dual Band Wireless-AC 8265 - 22.0.1.1";     Samples = 52096;  Pct = 100;                Vintage = "2020-09-28" }
    @{ Row = 24; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.60.2.1";          Samples = 26241;  Pct = 100;                Vintage = "2019-12-14" }
    @{ Row = 25; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.11.3";   Samples = 161874; Pct = 100;                Vintage = "2019-09-05" }
    @{ Row = 26; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.40.2.2";          Samples = 90508;  Pct = 99.90000000000001; Vintage = "2019-08-31" }
    @{ Row = 27; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5";   Samples = 154175; Pct = 99.90000000000001; Vintage = "2019-08-25" }
    @{ Row = 28; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.30.4.1";          Samples = 13016;  Pct = 100;                Vintage = "2019-07-29" }
    @{ Row = 29; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.10.2";   Samples = 20227;  Pct = 100;                Vintage = "2019-05-11" }
    @{ Row = 30; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.9.1";    Samples = 34065;  Pct = 100;                Vintage = "2019-04-28" }
    @{ Row = 31; Name = "Intel(R) Wi-Fi 6 AX200 160MHz - 21.10.1.2";          Samples = 52515;  Pct = 100;                Vintage = "2019-04-23" }
    @{ Row = 32; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.8.1";    Samples = 48540;  Pct = 100;                Vintage = "2019-03-16" }
    @{ Row = 33; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.70.5.2";    Samples = 184564; Pct = 99.90000000000001; Vintage = "2018-11-25" }
    @{ Row = 34; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.50.0.4";    Samples = 14221;  Pct = 100;                Vintage = "2018-05-08" }
    @{ Row = 35; Name = "Intel(R) Dual Band Wireless-AC 8265 - 20.30.1.2";    Samples = 23765;  Pct = 100;                Vintage = "2018-01-09" }
)

foreach ($d in $driverRows) {
    $r = $d.Row
    $needsFormat = ($r -gt 17)

    $ws.Cells.Item($r, 1).Value = $d.Name
    $ws.Cells.Item($r, 2).Value = $d.Samples
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = $d.Pct
    $ws.Cells.Item($r, 5).Value = $d.Vintage

    if ($needsFormat) {
        $ws.Cells.Item(13, 2).Copy() | Out-Null     # style 6 (#,##0 + right)
        $ws.Cells.Item($r, 2).PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Cells.Item(13, 4).Copy() | Out-Null     # style 4 (right)
        $ws.Cells.Item($r, 4).PasteSpecial($xlPasteFormats) | Out-Null
        $ws.Cells.Item(13, 5).Copy() | Out-Null     # style 4 (right)
        $ws.Cells.Item($r, 5).PasteSpecial($xlPasteFormats) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 6) Column A widens from 44 to 50 characters.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 49.17
